$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 182, pushing the existing rows 182-188 down to 183-189
$ws.Rows("182").Insert()

# Populate the newly inserted row 182 with the new weekly data point
$ws.Range("A182").Value = 4
$ws.Range("B182").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C182").Value = "Los Lagos"
$ws.Range("D182").Value = 44509
$ws.Range("D182").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E182").Value = 10
$ws.Range("F182").Value = 100112040
$ws.Range("G182").Value = "Cilantro"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 300
$ws.Range("K182").Value = 10000
$ws.Range("L182").Value = 10000
$ws.Range("M182").Value = 10000
$ws.Range("N182").Value = "$/caja 36 atados"
$ws.Range("O182").Value = "Región Metropolitana"
$ws.Range("P182").Value = 278
$ws.Range("Q182").Value = 36
$ws.Range("R182").Value = "Hortaliza"
